# Manually filled inflation rates
# - Row 37 (Inflation Rate / Spain-like row under "4. Prices"):
#     * Source note shortened from "OECD, Trading Economics, FX Empire"
#       to "OECD, Trading Economics"
#     * Data coverage (F37) filled in with value 99.8
#     * Remarks updated from "Key indicator, manually filled" to
#       "Key indicator, manually filled missing data"
# - Rows 16, 17, 25 (other "Key indicator, manually filled" remarks)
#     updated to "Key indicator, manually filled missing data"
# - Row 71 (Manufacturing PMI) data coverage "~65" -> "35/50"
# - Row 72 (Services PMI) data coverage "~30" -> "15/50"
# - Active selection moved to G70

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 37: shorten the source remark and fill in the missing data coverage value
$ws.Range("C37").Value = "OECD, Trading Economics"
$ws.Range("F37").Value = 99.8

# Update the "Key indicator, manually filled" remarks to add "missing data"
$ws.Range("G16").Value = "Key indicator, manually filled missing data"
$ws.Range("G17").Value = "Key indicator, manually filled missing data"
$ws.Range("G25").Value = "Key indicator, manually filled missing data"
$ws.Range("G37").Value = "Key indicator, manually filled missing data"

# Rows 71-72: replace the approximate data-coverage placeholders with counts
$ws.Range("F71").Value = "35/50"
$ws.Range("F72").Value = "15/50"

# Move/restore the active selection as recorded in the saved view state
$ws.Range("G70").Select()
